$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# Row 2: r_xy1y2
$ws.Range("G2").Value = -0.592095231879018
$ws.Range("H2").Value = -0.15648364867286
# Row 3: w_11
$ws.Range("G3").Value = -0.106600439582538
$ws.Range("H3").Value = 0.0196915373470342
# Row 4: w_21
$ws.Range("G4").Value = -0.0681066233993971
$ws.Range("H4").Value = 0.0673914372775766
# Row 5: r_xy1
$ws.Range("G5").Value = -0.383245794799145
$ws.Range("H5").Value = 0.0707942566741291
# Row 6: r_xy2
$ws.Range("G6").Value = -0.220657415569214
$ws.Range("H6").Value = 0.21834029700695
# Row 7: b_11
$ws.Range("G7").Value = -0.363325012333993
$ws.Range("H7").Value = 0.0671144329001289
# Row 8: b_21
$ws.Range("G8").Value = -0.23212699575646
$ws.Range("H8").Value = 0.229689435390395
# Row 9: main_effect
$ws.Range("G9").Value = -0.0861827736408375
$ws.Range("H9").Value = 0.0423707294621754
# Row 10: moderator_effect
$ws.Range("G10").Value = -0.0999089251591735
$ws.Range("H10").Value = -0.0494191756893807
# Row 11: interaction
$ws.Range("G11").Value = -0.0681756988538346
$ws.Range("H11").Value = -0.0180180172598485
# Row 12: q_b11_b21
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Style = "Normal"
# Row 13: q_rxy1_rxy2
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Style = "Normal"
# Row 14: cross_over_point
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Style = "Normal"
# Row 15: interaction_vs_main
$ws.Range("G15").Value = -0.052040737252145
$ws.Range("H15").Value = 0.094422409187166
# Row 16: interaction_vs_main_bscale
$ws.Range("G16").Value = -0.177369826785432
$ws.Range("H16").Value = 0.321818775953266
# Row 17: interaction_vs_main_rscale
$ws.Range("G17").Value = -0.151335397455756
$ws.Range("H17").Value = 0.304085488674867
# Row 18: dadas
$ws.Range("G18").Value = -0.136213246798794
$ws.Range("H18").Value = 0.134782874555153
# Row 19: dadas_bscale
$ws.Range("G19").Value = -0.46425399151292
$ws.Range("H19").Value = 0.459378870780789
# Row 20: dadas_rscale
$ws.Range("G20").Value = -0.441314831138428
$ws.Range("H20").Value = 0.436680594013899
# Row 21: abs_diff
$ws.Range("G21").Value = 0.0180180172598485
$ws.Range("H21").Value = 0.0681756988538346
# Row 22: abs_sum
$ws.Range("G22").Value = -0.0847414589243508
$ws.Range("H22").Value = 0.172365547281675
# Row 23: abs_diff_bscale
$ws.Range("G23").Value = 0.0614105942602595
$ws.Range("H23").Value = 0.232362424807539
# Row 24: abs_sum_bscale
$ws.Range("G24").Value = -0.288823308135155
$ws.Range("H24").Value = 0.587471447935084
# Row 25: abs_diff_rscale
$ws.Range("G25").Value = 0.0706048941135339
$ws.Range("H25").Value = 0.239529525449217
# Row 26: abs_sum_rscale
$ws.Range("G26").Value = -0.28113938964973
$ws.Range("H26").Value = 0.595908046337009
